$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- 1. Update "Status" text everywhere it appears (was "Ready for handoff") ---
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- 2. Update "Latest Handback DateTime" (column H) with real handback timestamps ---
$wsZhCn.Range("H2").Value = "2016-03-24 03:11:37"
$wsZhCn.Range("H3").Value = "2016-03-24 03:11:37"

$wsDeDe.Range("H2").Value = "2016-03-24 03:11:51"
$wsDeDe.Range("H3").Value = "2016-03-24 03:11:51"

# --- 3. Populate "Latest Target File" (F) and "Latest Handback File" (G) with hyperlinked file names ---

# zh-cn sheet
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a994505603472f3bfb8bf0118e82172d86f6ecfc/e2e/51e04dc1-69ed-4ea8-8fe0-b74347d37d56.md", "", "", "51e04dc1-69ed-4ea8-8fe0-b74347d37d56.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bedf2bbb38983ba49e5be1d9621b2af763d2fdbd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/51e04dc1-69ed-4ea8-8fe0-b74347d37d56.856bf71c15292fbd1b0ff7b1386f26a08ff67b3a.zh-cn.xlf", "", "", "51e04dc1-69ed-4ea8-8fe0-b74347d37d56.856bf71c15292fbd1b0ff7b1386f26a08ff67b3a.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a994505603472f3bfb8bf0118e82172d86f6ecfc/e2e/ea844430-8b1c-4d7d-9b6e-287556a92922.md", "", "", "ea844430-8b1c-4d7d-9b6e-287556a92922.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bedf2bbb38983ba49e5be1d9621b2af763d2fdbd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/ea844430-8b1c-4d7d-9b6e-287556a92922.2217f7e475f93555cb239d2c7c389af3b36c9ea7.zh-cn.xlf", "", "", "ea844430-8b1c-4d7d-9b6e-287556a92922.2217f7e475f93555cb239d2c7c389af3b36c9ea7.zh-cn.xlf")

# de-de sheet
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a994505603472f3bfb8bf0118e82172d86f6ecfc/e2e/51e04dc1-69ed-4ea8-8fe0-b74347d37d56.md", "", "", "51e04dc1-69ed-4ea8-8fe0-b74347d37d56.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f008f2d6c59f138eeb62ff96b43e000b7d519398/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/51e04dc1-69ed-4ea8-8fe0-b74347d37d56.856bf71c15292fbd1b0ff7b1386f26a08ff67b3a.de-de.xlf", "", "", "51e04dc1-69ed-4ea8-8fe0-b74347d37d56.856bf71c15292fbd1b0ff7b1386f26a08ff67b3a.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a994505603472f3bfb8bf0118e82172d86f6ecfc/e2e/ea844430-8b1c-4d7d-9b6e-287556a92922.md", "", "", "ea844430-8b1c-4d7d-9b6e-287556a92922.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f008f2d6c59f138eeb62ff96b43e000b7d519398/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/ea844430-8b1c-4d7d-9b6e-287556a92922.2217f7e475f93555cb239d2c7c389af3b36c9ea7.de-de.xlf", "", "", "ea844430-8b1c-4d7d-9b6e-287556a92922.2217f7e475f93555cb239d2c7c389af3b36c9ea7.de-de.xlf")
